# Scheduled market-data refresh: updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) for each leve row across all eight job sheets.
# Values below are the new market snapshot pulled by the runner; cells that
# have no new reading (e.g. a leve/material fell out of range) are cleared.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 1366
$ws.Range("J17").Value = 1366
$ws.Range("L17").Value = 4098
$ws.Range("N17").Value = -4434

# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 382.375
$ws.Range("I28").Value = 271.8
$ws.Range("J28").Value = 566.6667
$ws.Range("K28").Value = 271.8
$ws.Range("L28").Value = 566.6667
$ws.Range("M28").Value = 213.2
$ws.Range("N28").Value = -1536.6667

# Row 33: Glazed and Confused
$ws.Range("H33").Value = 248
$ws.Range("I33").Value = 185.36363
$ws.Range("K33").Value = 185.36363
$ws.Range("M33").Value = 43.63637

# Row 111: An Eye for Healing
$ws.Range("H111").Value = 808.375
$ws.Range("I111").Value = 761.2
$ws.Range("K111").Value = 2283.6
$ws.Range("M111").Value = 783.3999999999996

# Row 116: Growing Up
$ws.Range("H116").Value = 4654.4287
$ws.Range("I116").Value = 3649.25
$ws.Range("K116").Value = 3649.25
$ws.Range("M116").Value = -207.25

# Row 125: Body over Mind
$ws.Range("H125").Value = 105743.7
$ws.Range("I125").Value = 3623.5
$ws.Range("J125").Value = 258924
$ws.Range("K125").Value = 32611.5
$ws.Range("L125").Value = 2330316
$ws.Range("M125").Value = -30151.5
$ws.Range("N125").Value = -2335236

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3771
$ws.Range("I61").Value = 3966.1667
$ws.Range("J61").Value = 2600
$ws.Range("K61").Value = 3966.1667
$ws.Range("L61").Value = 2600
$ws.Range("M61").Value = -3754.1667
$ws.Range("N61").Value = -3024

# Row 97: Ore for Me
$ws.Range("H97").Value = 795.2727
$ws.Range("I97").Value = 795.2727
$ws.Range("K97").Value = 795.2727
$ws.Range("M97").Value = -299.2727

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 3771
$ws.Range("I136").Value = 3966.1667
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 11898.5001
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -9348.500100000001
$ws.Range("N136").Value = -12900

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 4009.6
$ws.Range("I99").Value = 4430.875
$ws.Range("K99").Value = 4430.875
$ws.Range("M99").Value = -2932.875

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0

# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 580.4
$ws.Range("I22").Value = 475.5
$ws.Range("K22").Value = 475.5
$ws.Range("M22").Value = -125.5

# Row 31: Wall Not Found
$ws.Range("H31").Value = 1879.4286
$ws.Range("I31").Value = 1723
$ws.Range("J31").Value = 3913
$ws.Range("K31").Value = 1723
$ws.Range("L31").Value = 3913
$ws.Range("M31").Value = -1428
$ws.Range("N31").Value = -4503

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 1879.4286
$ws.Range("I34").Value = 1723
$ws.Range("J34").Value = 3913
$ws.Range("K34").Value = 1723
$ws.Range("L34").Value = 3913
$ws.Range("M34").Value = -1521
$ws.Range("N34").Value = -4317

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 3049.4119
$ws.Range("I58").Value = 1716.1666
$ws.Range("J58").Value = 3776.6365
$ws.Range("K58").Value = 1716.1666
$ws.Range("L58").Value = 3776.6365
$ws.Range("M58").Value = -1513.1666
$ws.Range("N58").Value = -4182.636500000001

# Row 99: O Pine
$ws.Range("H99").Value = 4584.357
$ws.Range("I99").Value = 3727.1428
$ws.Range("J99").Value = 5441.5713
$ws.Range("K99").Value = 3727.1428
$ws.Range("L99").Value = 5441.5713
$ws.Range("M99").Value = -2229.1428
$ws.Range("N99").Value = -8437.5713

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1178
$ws.Range("I105").Value = 957.2727
$ws.Range("J105").Value = 1785
$ws.Range("K105").Value = 957.2727
$ws.Range("L105").Value = 1785
$ws.Range("M105").Value = 789.7273
$ws.Range("N105").Value = -5279

# Row 107: Built to Last
$ws.Range("H107").Value = 612.0476
$ws.Range("J107").Value = 724.8333
$ws.Range("L107").Value = 724.8333
$ws.Range("N107").Value = -4564.8333

# Row 113: Patient Patients
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 956.61536
$ws.Range("I122").Value = 956.61536
$ws.Range("K122").Value = 2869.84608
$ws.Range("M122").Value = -419.8460800000003

# Row 126: A Better Conductor
$ws.Range("H126").Value = 4584.357
$ws.Range("I126").Value = 3727.1428
$ws.Range("J126").Value = 5441.5713
$ws.Range("K126").Value = 11181.4284
$ws.Range("L126").Value = 16324.7139
$ws.Range("M126").Value = -8711.428400000001
$ws.Range("N126").Value = -21264.7139

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1658.7142
$ws.Range("I134").Value = 1562.2
$ws.Range("K134").Value = 4686.6
$ws.Range("M134").Value = -2151.6

# Row 136: Turali Quality
$ws.Range("H136").Value = 3049.4119
$ws.Range("I136").Value = 1716.1666
$ws.Range("J136").Value = 3776.6365
$ws.Range("K136").Value = 5148.4998
$ws.Range("L136").Value = 11329.9095
$ws.Range("M136").Value = -2598.4998
$ws.Range("N136").Value = -16429.9095

$ws = $wb.Worksheets.Item("CUL")
# Row 17: Chew the Fat
$ws.Range("H17").Value = 227.8
$ws.Range("J17").Value = 282.25
$ws.Range("L17").Value = 846.75
$ws.Range("N17").Value = -1184.75

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 3259.6
$ws.Range("I80").Value = 2933
$ws.Range("K80").Value = 2933
$ws.Range("M80").Value = -1935

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 3259.6
$ws.Range("I83").Value = 2933
$ws.Range("K83").Value = 14665
$ws.Range("M83").Value = -9673

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 850
$ws.Range("I97").Value = 850
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 850
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -354

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 14022.2
$ws.Range("J126").Value = 16500
$ws.Range("L126").Value = 49500
$ws.Range("N126").Value = -54440

$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head
$ws.Range("H2").Value = 1.4285715
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2.5
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 2.5
$ws.Range("M2").Value = 111
$ws.Range("N2").Value = -226.5

# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 20025.615
$ws.Range("I7").Value = 17093.111
$ws.Range("J7").Value = 26623.75
$ws.Range("K7").Value = 17093.111
$ws.Range("L7").Value = 26623.75
$ws.Range("M7").Value = -16981.111
$ws.Range("N7").Value = -26847.75

# Row 16: Saddle Sore
$ws.Range("H16").Value = 2557.1428
$ws.Range("I16").Value = 3000
$ws.Range("J16").Value = 2483.3333
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 2483.3333
$ws.Range("M16").Value = -2830
$ws.Range("N16").Value = -2823.3333

# Row 40: Best Served Toad
$ws.Range("H40").Value = 7505.25
$ws.Range("I40").Value = 7505.25
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 7505.25
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -7369.25

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 212.27272
$ws.Range("I55").Value = 170.14285
$ws.Range("J55").Value = 286
$ws.Range("K55").Value = 170.14285
$ws.Range("L55").Value = 286
$ws.Range("M55").Value = 2.85714999999999
$ws.Range("N55").Value = -632

# Row 74: Overall, We Blend In
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

# Row 77: Eviction Notice (L)
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

# Row 98: Try Tricorne Again
$ws.Range("H98").Value = 28999
$ws.Range("J98").Value = 28999
$ws.Range("L98").Value = 28999
$ws.Range("N98").Value = -34989

# Row 122: Hell on Leather
$ws.Range("H122").Value = 3438.5
$ws.Range("I122").Value = 3287.5715
$ws.Range("J122").Value = 4495
$ws.Range("K122").Value = 9862.7145
$ws.Range("L122").Value = 13485
$ws.Range("M122").Value = -7412.7145
$ws.Range("N122").Value = -18385

# Row 126: Battered Books
$ws.Range("H126").Value = 20025.615
$ws.Range("I126").Value = 17093.111
$ws.Range("J126").Value = 26623.75
$ws.Range("K126").Value = 51279.333
$ws.Range("L126").Value = 79871.25
$ws.Range("M126").Value = -48809.333
$ws.Range("N126").Value = -84811.25

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 6030.4546
$ws.Range("I132").Value = 6284.222
$ws.Range("K132").Value = 18852.666
$ws.Range("M132").Value = -16322.666

$ws = $wb.Worksheets.Item("WVR")
# Row 5: Hire in the Blood
$ws.Range("H5").Value = 52000
$ws.Range("J5").Value = 52000
$ws.Range("L5").Value = 52000
$ws.Range("N5").Value = -52224

# Row 26: New Shoes, New Me
$ws.Range("H26").Value = 19008.8
$ws.Range("J26").Value = 23758
$ws.Range("L26").Value = 23758
$ws.Range("N26").Value = -24344

# Row 100: Of Great Import
$ws.Range("H100").Value = 1117
$ws.Range("I100").Value = 800.5
$ws.Range("J100").Value = 1750
$ws.Range("K100").Value = 1601
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -1060
$ws.Range("N100").Value = -4582
